# Lattice-multiplication exercise sheet: refresh the 5x3 table of
# practice problems with a new set of factor pairs / lattice digits.
#
# Each table cell holds a single run whose text is split by <w:br/>
# line breaks into five lines:
#   1) "NN x MM"      - the two factors
#   2) "  d    d"     - the ones-digit of each factor, spaced out
#   3) "  ----"       - divider (unchanged by this edit)
#   4) "d|    |"      - lattice box, row 1
#   5) "d|    |"      - lattice box, row 2
#
# We rewrite each cell's Range.Text in one shot (old lines joined with
# a vertical-tab, which Word's Range.Text uses to represent <w:br/>)
# so the whole run is regenerated consistently from the new values.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vtab = [char]11

# row, col, new "NN x MM", new "  d    d", new line1 "d|    |", new line2 "d|    |"
$updates = @(
    @(1, 1, "74 x 66", "  6    6", "7|    |", "4|    |"),
    @(1, 2, "60 x 30", "  3    0", "6|    |", "0|    |"),
    @(1, 3, "41 x 19", "  1    9", "4|    |", "1|    |"),
    @(2, 1, "84 x 61", "  6    1", "8|    |", "4|    |"),
    @(2, 2, "24 x 54", "  5    4", "2|    |", "4|    |"),
    @(2, 3, "39 x 54", "  5    4", "3|    |", "9|    |"),
    @(3, 1, "25 x 60", "  6    0", "2|    |", "5|    |"),
    @(3, 2, "51 x 16", "  1    6", "5|    |", "1|    |"),
    @(3, 3, "21 x 85", "  8    5", "2|    |", "1|    |"),
    @(4, 1, "53 x 45", "  4    5", "5|    |", "3|    |"),
    @(4, 2, "73 x 29", "  2    9", "7|    |", "3|    |"),
    @(4, 3, "67 x 97", "  9    7", "6|    |", "7|    |"),
    @(5, 1, "47 x 76", "  7    6", "4|    |", "7|    |"),
    @(5, 2, "68 x 97", "  9    7", "6|    |", "8|    |"),
    @(5, 3, "86 x 61", "  6    1", "8|    |", "6|    |")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $header = $u[2]
    $digits = $u[3]
    $line1 = $u[4]
    $line2 = $u[5]

    $cell = $t.Cell($row, $col)
    $newText = $header + $vtab + $digits + $vtab + "  ----" + $vtab + $line1 + $vtab + $line2
    $cell.Range.Text = $newText
}
